# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The sheet currently runs from A1:AC70 (player bio/stat columns). We append
# three new columns - AD (Wins), AE (Losses), AF (Ties) - with a bold/boxed
# header in row 1 (matching the existing header styling) and the team's
# 2022 season record (62-100-0) repeated down every data row (2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the styling (bold, centered, boxed) already used by the rest of the
# header row by copying the format from the neighboring header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (rows 2-70) ---------------------------------------------
$lastRow = 70
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 62   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 100  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
